$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4370.8335
$ws.Range("I64").Value = 4427.143
$ws.Range("K64").Value = 4427.143
$ws.Range("M64").Value = -4179.143

$ws.Range("H67").Value = 4370.8335
$ws.Range("I67").Value = 4427.143
$ws.Range("K67").Value = 4427.143
$ws.Range("M67").Value = -3569.143

$ws.Range("H98").Value = 3997.2273
$ws.Range("I98").Value = 4206.05
$ws.Range("K98").Value = 4206.05
$ws.Range("M98").Value = -2708.05

$ws.Range("H112").Value = 5484.3794
$ws.Range("J112").Value = 6236.08
$ws.Range("L112").Value = 18708.24
$ws.Range("N112").Value = -20924.24

$ws.Range("H122").Value = 3997.2273
$ws.Range("I122").Value = 4206.05
$ws.Range("K122").Value = 12618.15
$ws.Range("M122").Value = -10168.15

$ws.Range("H137").Value = 1092.7792
$ws.Range("I137").Value = 798.82355
$ws.Range("J137").Value = 1669.3846
$ws.Range("K137").Value = 2396.47065
$ws.Range("L137").Value = 5008.1538
$ws.Range("M137").Value = 153.5293500000002
$ws.Range("N137").Value = -10108.1538

$ws.Range("H141").Value = 741.6111
$ws.Range("I141").Value = 608.7646999999999
$ws.Range("K141").Value = 1826.2941
$ws.Range("M141").Value = 3353.7059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 34533.332
$ws.Range("I2").Value = 1792
$ws.Range("J2").Value = 50904
$ws.Range("K2").Value = 1792
$ws.Range("L2").Value = 50904
$ws.Range("M2").Value = -1679
$ws.Range("N2").Value = -51130

$ws.Range("H32").Value = 3801.2405
$ws.Range("I32").Value = 3442.261
$ws.Range("K32").Value = 3442.261
$ws.Range("M32").Value = -3155.261

$ws.Range("H61").Value = 18519630
$ws.Range("I61").Value = 22223182
$ws.Range("K61").Value = 22223182
$ws.Range("M61").Value = -22222970

$ws.Range("H116").Value = 34533.332
$ws.Range("I116").Value = 1792
$ws.Range("J116").Value = 50904
$ws.Range("K116").Value = 1792
$ws.Range("L116").Value = 50904
$ws.Range("M116").Value = 502
$ws.Range("N116").Value = -55492

$ws.Range("H117").Value = 52500
$ws.Range("J117").Value = 52500
$ws.Range("L117").Value = 52500
$ws.Range("N117").Value = -61678

$ws.Range("H136").Value = 18519630
$ws.Range("I136").Value = 22223182
$ws.Range("K136").Value = 66669546
$ws.Range("M136").Value = -66666996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 34533.332
$ws.Range("I3").Value = 1792
$ws.Range("J3").Value = 50904
$ws.Range("K3").Value = 1792
$ws.Range("L3").Value = 50904
$ws.Range("M3").Value = -1678
$ws.Range("N3").Value = -51132

$ws.Range("H96").Value = 7800
$ws.Range("I96").Value = 4200
$ws.Range("J96").Value = 15000
$ws.Range("K96").Value = 4200
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -1454
$ws.Range("N96").Value = -20492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2322.2273
$ws.Range("I31").Value = 2530.5625
$ws.Range("K31").Value = 2530.5625
$ws.Range("M31").Value = -2235.5625

$ws.Range("H34").Value = 2322.2273
$ws.Range("I34").Value = 2530.5625
$ws.Range("K34").Value = 2530.5625
$ws.Range("M34").Value = -2328.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5636.4736
$ws.Range("I107").Value = 426.4
$ws.Range("J107").Value = 7497.2144
$ws.Range("K107").Value = 1279.2
$ws.Range("L107").Value = 22491.6432
$ws.Range("M107").Value = 640.8000000000002
$ws.Range("N107").Value = -26331.6432

$ws.Range("H131").Value = 16950362
$ws.Range("J131").Value = 1304.2885
$ws.Range("L131").Value = 3912.8655
$ws.Range("N131").Value = -13992.8655

$ws.Range("H136").Value = 2551.125
$ws.Range("I136").Value = 1302.5
$ws.Range("K136").Value = 3907.5
$ws.Range("M136").Value = 1192.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1229.5807
$ws.Range("I102").Value = 1089.5186
$ws.Range("K102").Value = 1089.5186
$ws.Range("M102").Value = 532.4813999999999

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H140").Value = 34237.5
$ws.Range("J140").Value = 34237.5
$ws.Range("L140").Value = 34237.5
$ws.Range("N140").Value = -44597.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 654.0454999999999
$ws.Range("I22").Value = 457
$ws.Range("J22").Value = 938.6667
$ws.Range("K22").Value = 457
$ws.Range("L22").Value = 938.6667
$ws.Range("M22").Value = -162
$ws.Range("N22").Value = -1528.6667

$ws.Range("H27").Value = 654.0454999999999
$ws.Range("I27").Value = 457
$ws.Range("J27").Value = 938.6667
$ws.Range("K27").Value = 457
$ws.Range("L27").Value = 938.6667
$ws.Range("M27").Value = -350
$ws.Range("N27").Value = -1152.6667

$ws.Range("H55").Value = 326.9355
$ws.Range("I55").Value = 219.55556
$ws.Range("J55").Value = 475.6154
$ws.Range("K55").Value = 219.55556
$ws.Range("L55").Value = 475.6154
$ws.Range("M55").Value = -46.55556000000001
$ws.Range("N55").Value = -821.6154

$ws.Range("H63").Value = 19000
$ws.Range("I63").Value = 19000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 19000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -18251
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 19000
$ws.Range("I66").Value = 19000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 57000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -53256
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 20810.25
$ws.Range("I132").Value = 1214.9714
$ws.Range("K132").Value = 3644.9142
$ws.Range("M132").Value = -1114.9142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 177253.25
$ws.Range("I14").Value = 235337.67
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 235337.67
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -235169.67
$ws.Range("N14").Value = -3336

$ws.Range("H46").Value = 34995
$ws.Range("J46").Value = 34995
$ws.Range("L46").Value = 34995
$ws.Range("N46").Value = -35457

$ws.Range("H63").Value = 13307
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 13307
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 13307
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -14555

$ws.Range("H66").Value = 13307
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 13307
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 39921
$ws.Range("N66").Value = -46161

$ws.Range("H132").Value = 2797.3438
$ws.Range("I132").Value = 3097.2964
$ws.Range("J132").Value = 1177.6
$ws.Range("K132").Value = 9291.889200000001
$ws.Range("L132").Value = 3532.8
$ws.Range("M132").Value = -6761.889200000001
$ws.Range("N132").Value = -8592.799999999999

$ws.Range("H134").Value = 34995
$ws.Range("J134").Value = 34995
$ws.Range("L134").Value = 104985
$ws.Range("N134").Value = -110055

$ws.Range("H136").Value = 736.4545000000001
$ws.Range("I136").Value = 510.8889
$ws.Range("J136").Value = 1751.5
$ws.Range("K136").Value = 1532.6667
$ws.Range("L136").Value = 5254.5
$ws.Range("M136").Value = 1017.3333
$ws.Range("N136").Value = -10354.5
